$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "ķųųųų"
$ws.Range("A4").Value = "   long long something   "

$ws.Range("A4").Select()
